# Fruta / hortaliza, semanal
#
# A new weekly price record for "Caqui" (Mankaki, Primera) needs to be
# inserted as row 29 of the data table, pushing the existing rows 29-54
# down to rows 30-55 (dimension grows from A1:T54 to A1:T55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29; this shifts rows 29..54 down to 30..55
# and inherits the formatting (incl. the date style on column D) from the
# row above, matching the rest of the sheet.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record.
$ws.Range("A29").Value = 9
$ws.Range("B29").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 45049
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = "Otros"
$ws.Range("I29").Value = 100107001
$ws.Range("J29").Value = "Caqui"
$ws.Range("K29").Value = "Mankaki"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 500
$ws.Range("N29").Value = 13000
$ws.Range("O29").Value = 14000
$ws.Range("P29").Value = 13560
$ws.Range("Q29").Value = "`$/caja 16 kilos granel"
$ws.Range("R29").Value = "Región de O'Higgins"
$ws.Range("S29").Value = 848
$ws.Range("T29").Value = 16
